$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 314
$ws1.Range("F3").Value = 1258
$ws1.Range("F4").Value = 362
$ws1.Range("F5").Value = 330
$ws1.Range("F6").Value = 3840
$ws1.Range("F8").Value = 752
$ws1.Range("F9").Value = 2220
$ws1.Range("F13").Value = 159
$ws1.Range("F15").Value = 2125
$ws1.Range("F19").Value = 334
$ws1.Range("F20").Value = 224
$ws1.Range("F21").Value = 21

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 10
$ws2.Range("F9").Value = 93
$ws2.Range("F12").Value = 225

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 817
$ws3.Range("F4").Value = 2080
$ws3.Range("F5").Value = 313

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 817
$ws4.Range("F4").Value = 2080
$ws4.Range("F5").Value = 313
$ws4.Range("F10").Value = 314
$ws4.Range("F11").Value = 1258
$ws4.Range("F12").Value = 362
$ws4.Range("F15").Value = 10
$ws4.Range("F16").Value = 330
$ws4.Range("F17").Value = 3840
$ws4.Range("F20").Value = 93
$ws4.Range("F23").Value = 752
$ws4.Range("F24").Value = 2220
$ws4.Range("F26").Value = 225
$ws4.Range("F29").Value = 159
$ws4.Range("F32").Value = 2125
$ws4.Range("F38").Value = 334
$ws4.Range("F39").Value = 224
$ws4.Range("F40").Value = 21

$wb.Save()
